$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("DeviceSetupLogins")
$ws3 = $wb.Worksheets.Item("InjectSpecificUser")

# --- Sheet1: flip "ignore" column from yes -> no for the Auto.* rows ---
$ws1.Range("C2:C4").Value = "no"

# --- DeviceSetupLogins: bump password + flip ignore for the Auto.* rows ---
$ws2.Range("B2:B4").Value = "MHRA12345"
$ws2.Range("C2:C4").Value = "no"

# --- InjectSpecificUser: bump password for the Auto.* rows ---
$ws3.Range("B2:B4").Value = "MHRA12345"

# --- Update each sheet's selection/active cell ---
$ws1.Range("B2").Select() | Out-Null
$ws2.Range("B2:B4").Select() | Out-Null
$ws3.Range("B2:B4").Select() | Out-Null

# --- Make InjectSpecificUser the active (selected) sheet/tab ---
$ws3.Activate() | Out-Null
